$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" (first sheet): update Date, insert "Jurisdiction" row ---
$ws = $wb.Worksheets.Item(1)

# 1. Update the Date value (row 8, column B)
$ws.Cells.Item(8, 2).Value = "2024-09-12T14:01:50+00:00"

# 2. Insert a new "Jurisdiction" row right after "Contact" (row 10), i.e. at row 11,
#    pushing every existing row from 11..20 down to 12..21.
#    We do this by shifting cell VALUES (not using Rows.Insert, which would create
#    a freshly-styled row) so that the existing "s=2" cell style is preserved.

# Capture current values for rows 11..20 (A and B columns) before overwriting anything.
$colA = @{}
$colB = @{}
for ($r = 11; $r -le 20; $r++) {
    $colA[$r] = $ws.Cells.Item($r, 1).Value2
    $colB[$r] = $ws.Cells.Item($r, 2).Value2
}

# Write them back shifted down by one row, starting from the bottom so we never
# clobber a value before it has been read.
for ($r = 20; $r -ge 11; $r--) {
    $ws.Cells.Item($r + 1, 1).Value = $colA[$r]
    $ws.Cells.Item($r + 1, 2).Value = $colB[$r]
}

# Now populate the newly freed-up row 11 with the Jurisdiction property.
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""
